$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record at row 255 (pushes old rows 255-272 down
# to 256-273, growing the used range from A1:R272 to A1:R273).
$ws.Rows.Item(255).Insert()

$ws.Range("A255").Value = 3
$ws.Range("B255").Value = "Femacal de La Calera"
$ws.Range("C255").Value = "Coquimbo"
$ws.Range("D255").Value = 45008
$ws.Range("E255").Value = 5
$ws.Range("F255").Value = 100112030
$ws.Range("G255").Value = "Poroto granado"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 73
$ws.Range("K255").Value = 32000
$ws.Range("L255").Value = 33000
$ws.Range("M255").Value = 32479
$ws.Range("N255").Value = '$/saco 25 kilos'
$ws.Range("O255").Value = "Provincia de Quillota"
$ws.Range("P255").Value = 1299
$ws.Range("Q255").Value = 25
$ws.Range("R255").Value = "Hortaliza"
